$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.600.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.127.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "527.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.13%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.130.56"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.95%  "
$ws.Range("E9").Value = "  +2.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.107"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("E12").Value = "  +3.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.670.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("E14").Value = "  +2.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000164"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "57.791.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.145.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "349.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.506"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.05%  "
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.985"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0914"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.10"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "21.10"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.49%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.99"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.18"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.43%  "
$ws.Range("E41").Value = "  +1.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.53%  "
$ws.Range("E43").Value = "  +7.44%  "
$ws.Range("E44").Value = "  +2.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.168.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "36.60"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.54%  "
$ws.Range("E47").Value = "  +4.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.347.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.12%  "
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.967"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("E51").Value = "  +0.51%  "
